$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'318.01"
$ws.Range("E2").Value = "'4.13%"

$ws.Range("E3").Value = "'1.95%"

$ws.Range("D4").Value = "'5.135"
$ws.Range("E4").Value = "'0.54%"

$ws.Range("D5").Value = "'0.08201"
$ws.Range("E5").Value = "'1.48%"

$ws.Range("D6").Value = "'2.038"
$ws.Range("E6").Value = "'6.07%"

$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.353"
$ws.Range("E7").Value = "'4.06%"

$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9375"
$ws.Range("E8").Value = "'1.03%"

$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1354"
$ws.Range("E9").Value = "'-6.96%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1986"
$ws.Range("E10").Value = "'3.67%"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09112"
$ws.Range("E11").Value = "'1.00%"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03497"
$ws.Range("E12").Value = "'-0.43%"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09773"
$ws.Range("E13").Value = "'0.00%"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001405"
$ws.Range("E14").Value = "'0.70%"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006111"
$ws.Range("E15").Value = "'4.85%"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.687"
$ws.Range("E16").Value = "'-2.53%"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.337"
$ws.Range("E17").Value = "'3.20%"

$ws.Range("D19").Value = "'0.3475"
$ws.Range("E19").Value = "'1.31%"

$ws.Range("D20").Value = "'0.1314"
$ws.Range("E20").Value = "'-1.02%"

$ws.Range("D21").Value = "'4.957"
$ws.Range("E21").Value = "'5.38%"

$ws.Range("D22").Value = "'0.2449"
$ws.Range("E22").Value = "'1.25%"

$ws.Range("D23").Value = "'0.04364"
$ws.Range("E23").Value = "'-0.41%"

$ws.Range("D24").Value = "'0.001229"
$ws.Range("E24").Value = "'-0.77%"

$ws.Range("D25").Value = "'0.004807"
$ws.Range("E25").Value = "'12.54%"

$ws.Range("E26").Value = "'-0.16%"

$ws.Range("D27").Value = "'0.0003996"
$ws.Range("E27").Value = "'-10.15%"

$ws.Range("D39").Value = "'0.02256"
$ws.Range("E39").Value = "'11.46%"

$ws.Range("D40").Value = "'0.05204"
$ws.Range("E40").Value = "'3.11%"

$ws.Range("D41").Value = "'0.007763"
$ws.Range("E41").Value = "'3.06%"

$ws.Range("D42").Value = "'0.009882"
$ws.Range("E42").Value = "'0.91%"

$ws.Range("D43").Value = "'0.1403"
$ws.Range("E43").Value = "'4.60%"

$ws.Range("D44").Value = "'0.002045"
$ws.Range("E44").Value = "'-2.92%"

$ws.Range("D45").Value = "'0.009102"
$ws.Range("E45").Value = "'-8.45%"

$ws.Range("D46").Value = "'0.00006576"
$ws.Range("E46").Value = "'6.14%"

$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.37%"

$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.002944"
$ws.Range("E48").Value = "'2.34%"

$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.001688"
$ws.Range("E49").Value = "'-6.46%"

$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.37%"

$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.37%"

